$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Ampli audio 3W 8ohms ---
$amazonUrl = "https://www.amazon.fr/Amplificateur-2-5V-5-5V-dAmplificateur-Puissance-Amplifier/dp/B08D9MGHM2/ref=sxin_15_pa_sp_search_thematic_sspa?cv_ct_cx=ampli+8+ohms&sbo=RZvfv%2F%2FHxDF%2BO5021pAnSA%3D%3D&sr=1-1-86ee67e3-2ea6-4725-8419-71cfe38eb657-spons&sp_csd=d2lkZ2V0TmFtZT1zcF9zZWFyY2hfdGhlbWF0aWM&psc=1"
$ws.Range("C4").Value = $amazonUrl
$ws.Range("A4").Value = "Ampli audio 3W 8ohms"

$ws.Range("B2").Copy($ws.Range("B4"))
$ws.Range("B4").Value = 6.31

$ws.Hyperlinks.Add($ws.Range("C4"), $amazonUrl, [Type]::Missing, [Type]::Missing, $amazonUrl)
$ws.Range("C4").Style = "Hyperlink"

# --- Row 5: haut parleur 3w  8ohms ---
$mouserUrl = "https://www.mouser.fr/ProductDetail/Same-Sky/CMS-3118-38E?qs=IKkN%2F947nfB5KUQR4YXCyg%3D%3D"
$ws.Range("C5").Value = $mouserUrl
$ws.Range("A5").Value = "haut parleur 3w  8ohms"

$ws.Range("B2").Copy($ws.Range("B5"))
$ws.Range("B5").Value = 5.91

$ws.Hyperlinks.Add($ws.Range("C5"), $mouserUrl)
$ws.Range("C5").Style = "Hyperlink"

# --- Selection ---
$ws.Range("C14").Select() | Out-Null
